# Rename the worksheet from "b" to "AABC"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "AABC"

# Extend the table with two new rows of data (rows 7 and 8),
# carrying the same formatting as the existing "A" column entries above them.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("D8").Value = 1.7
